$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.146.65'
$ws.Range("E2").Value = '  +2.84%  '

$ws.Range("D3").Value = '3.756.21'
$ws.Range("E3").Value = '  +2.51%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.94'
$ws.Range("E5").Value = '  +1.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.82'
$ws.Range("E6").Value = '  +2.76%  '

$ws.Range("D7").Value = '3.759.01'
$ws.Range("E7").Value = '  +2.68%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.536'
$ws.Range("E9").Value = '  +2.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.168'
$ws.Range("E10").Value = '  +6.29%  '

$ws.Range("E11").Value = '  +3.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.464'
$ws.Range("E12").Value = '  +0.78%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.46'
$ws.Range("E13").Value = '  +2.97%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000250'
$ws.Range("E14").Value = '  +4.54%  '

$ws.Range("D15").Value = '4.376.32'
$ws.Range("E15").Value = '  +2.52%  '

$ws.Range("D16").Value = '3.754.90'
$ws.Range("E16").Value = '  +2.78%  '

$ws.Range("D17").Value = '69.091.45'
$ws.Range("E17").Value = '  +2.81%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.30'
$ws.Range("E18").Value = '  +2.07%  '

$ws.Range("E19").Value = '  +0.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.13'
$ws.Range("E20").Value = '  -2.50%  '

$ws.Range("E21").Value = '  +19.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '494.97'
$ws.Range("E22").Value = '  +0.56%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.729'
$ws.Range("E23").Value = '  +1.84%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000157'
$ws.Range("E24").Value = '  +15.73%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.44'
$ws.Range("E25").Value = '  +0.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.34'
$ws.Range("E26").Value = '  +2.20%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.40'
$ws.Range("E27").Value = '  +2.31%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.44'
$ws.Range("E28").Value = '  +5.06%  '

$ws.Range("E29").Value = '  +0.48%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.55'
$ws.Range("E30").Value = '  +8.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.99'
$ws.Range("E31").Value = '  +2.63%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.99'
$ws.Range("E32").Value = '  +4.56%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '32.09'
$ws.Range("E33").Value = '  +1.92%  '

$ws.Range("D34").Value = '3.894.69'
$ws.Range("E34").Value = '  +2.58%  '

$ws.Range("E35").Value = '  +1.85%  '

$ws.Range("D36").Value = '3.684.27'
$ws.Range("E36").Value = '  +2.50%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.10%  '

$ws.Range("E38").Value = '  +2.69%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.89'
$ws.Range("E39").Value = '  +2.52%  '

$ws.Range("E40").Value = '  +1.98%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.325'
$ws.Range("E41").Value = '  +1.10%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.07'
$ws.Range("E42").Value = '  +11.49%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '438.82'
$ws.Range("E43").Value = '  +1.67%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '48.87'
$ws.Range("E44").Value = '  +0.44%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.98'
$ws.Range("E45").Value = '  +3.29%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.48'
$ws.Range("E46").Value = '  +1.76%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.69'
$ws.Range("E48").Value = '  +0.80%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.77'
$ws.Range("E49").Value = '  -0.39%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0355'
$ws.Range("E50").Value = '  +2.98%  '

$ws.Range("D51").Value = '2.782.66'
$ws.Range("E51").Value = '  +1.44%  '
